$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate row-1 header labels (English -> German) ---------------------
$ws.Range("A1").Value = "Interner Aufbau ab"
$ws.Range("B1").Value = "Externer Aufbau ab"
$ws.Range("C1").Value = "VA ab"
$ws.Range("D1").Value = "VA bis"
$ws.Range("E1").Value = "Externer Abbau bis"
$ws.Range("F1").Value = "Inerner Abbau bis"
$ws.Range("G1").Value = "Matchcode"
$ws.Range("H1").Value = "Titel"
$ws.Range("I1").Value = "Kommentar"
$ws.Range("J1").Value = "Konto"
$ws.Range("K1").Value = "Typ"
$ws.Range("L1").Value = "Status"
$ws.Range("M1").Value = "Ort"
$ws.Range("N1").Value = "Projekt"
$ws.Range("O1").Value = "Technik"
$ws.Range("P1").Value = "TPL"
$ws.Range("Q1").Value = "PLM"
$ws.Range("R1").Value = "Sicherheit"
$ws.Range("S1").Value = "Projekt2"
$ws.Range("T1").Value = "TPL2"
$ws.Range("U1").Value = "Technik2"
$ws.Range("V1").Value = "PLM2"
$ws.Range("W1").Value = "Sicherheit2"

# --- Apply a date number format to the first six template cells of row 2 ---
# (these hold the jsDate -> xlDate placeholders: internal/external build-up,
# event-from/to, external/internal dismantling)
$ws.Range("A2:F2").NumberFormat = "dd/mm/yyyy"

# --- Turn on AutoFilter for the header + template row ----------------------
[void]$ws.Range("A1:W2").AutoFilter()

# Register the (hidden, sheet-scoped) _FilterDatabase defined name that
# Excel/Calc persist alongside an AutoFilter range.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=events!`$A`$1:`$W`$2")
$fdb.Visible = $false

# --- Move the active selection back to A1 -----------------------------------
[void]$ws.Range("A1").Select()
